$wb = $excel.ActiveWorkbook

# --- Data sheet: add "Hotdogs" and "Shirt" columns (D, E) ---
$wsData = $wb.Worksheets.Item("Data")

$wsData.Range("D1").Value = "Hotdogs"
$wsData.Range("E1").Value = "Shirt"
$wsData.Range("D1:E1").Font.Bold = $true

$hotdogs = @(12, 3, 21, 6, 15, 9, 25, 4, 17, 4, 27, 3, 18, 29)
$shirts  = @("M", "L", "S", "M", "L", "M", "XL", "L", "M", "XXL", "L", "M", "S", "L")

for ($i = 0; $i -lt $hotdogs.Length; $i++) {
    $row = $i + 2
    $wsData.Cells.Item($row, 4).Value = $hotdogs[$i]
    $wsData.Cells.Item($row, 5).Value = $shirts[$i]
}

# --- Codebook sheet: document the two new variables ---
$wsCode = $wb.Worksheets.Item("Codebook")

$wsCode.Range("A5").Value = "Hotdogs"
$wsCode.Range("B5").Value = "number of hotdogs they could eat"
$wsCode.Range("C5").Value = "numeric value >0 or NA"

$wsCode.Range("A6").Value = "Shirt"
$wsCode.Range("B6").Value = "standard shirt sizes"
$wsCode.Range("C6").Value = "S/M/L/XL/XXL"

# --- Selections / active sheet state ---
$wsData.Range("G6").Select()
$wsCode.Range("C8").Select()
$wsCode.Activate()
